$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.409481333333334
$ws.Range("H2").Value = 28.228444
$ws.Range("I2").Value = 0.2433300530093958
$ws.Range("J2").Value = 0.2433300530093958
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 14.04343760842311
$ws.Range("R2").Value = 126.390938475808
$ws.Range("S2").Value = 0.006738987327018824
$ws.Range("T2").Value = 0.006738987327018825
$ws.Range("G3").Value = 9.409481333333334
$ws.Range("H3").Value = 28.228444
$ws.Range("I3").Value = 0.2433300530093958
$ws.Range("J3").Value = 0.2433300530093958
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 318.5849142025885
$ws.Range("R3").Value = 2867.264227823297
$ws.Range("S3").Value = 0.1528785016357327
$ws.Range("T3").Value = 0.1528785016357327
$ws.Range("G4").Value = 9.409481333333334
$ws.Range("H4").Value = 28.228444
$ws.Range("I4").Value = 0.2433300530093958
$ws.Range("J4").Value = 0.2433300530093958
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 174.4493813657662
$ws.Range("R4").Value = 1570.044432291896
$ws.Range("S4").Value = 0.08371256404664426
$ws.Range("T4").Value = 0.08371256404664427
$ws.Range("I5").Value = 0.5069354697952918
$ws.Range("J5").Value = 0.5069354697952919
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 29.25703814025778
$ws.Range("R5").Value = 263.31334326232
$ws.Range("S5").Value = 0.01403949764657674
$ws.Range("T5").Value = 0.01403949764657674
$ws.Range("I6").Value = 0.5069354697952918
$ws.Range("J6").Value = 0.5069354697952919
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("S6").Value = 0.3184955334938341
$ws.Range("T6").Value = 0.3184955334938342
$ws.Range("I7").Value = 0.5069354697952918
$ws.Range("J7").Value = 0.5069354697952919
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("S7").Value = 0.1744004386548809
$ws.Range("T7").Value = 0.174400438654881
$ws.Range("G8").Value = 9.657138
$ws.Range("I8").Value = 0.2497344771953123
$ws.Range("J8").Value = 0.2497344771953124
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 14.413059569872
$ws.Range("R8").Value = 129.717536128848
$ws.Range("S8").Value = 0.00691635684176626
$ws.Range("T8").Value = 0.006916356841766261
$ws.Range("G9").Value = 9.657138
$ws.Range("I9").Value = 0.2497344771953123
$ws.Range("J9").Value = 0.2497344771953124
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("R9").Value = 2942.730353527776
$ws.Range("S9").Value = 0.1569022494682488
$ws.Range("T9").Value = 0.1569022494682488
$ws.Range("G10").Value = 9.657138
$ws.Range("I10").Value = 0.2497344771953123
$ws.Range("J10").Value = 0.2497344771953124
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("S10").Value = 0.08591587088529733
$ws.Range("T10").Value = 0.08591587088529734
